# Generate Report for handback
# Adds a new handed-back file (eb931bac-a78a-4017-93c9-acbdf45de9c6) as row 4
# to the Overview sheet as well as the zh-cn and de-de detail sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Overview" -> new row 4 (File Name | zh-cn | de-de)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/87e2101526f2d38fc04e67819f1fc856141ab27e/e2e/eb931bac-a78a-4017-93c9-acbdf45de9c6.md",
    "",
    "",
    "eb931bac-a78a-4017-93c9-acbdf45de9c6.md"
) | Out-Null

$wsOverview.Range("B4").Value2 = "Handed back: in sync with en-US"
$wsOverview.Range("C4").Value2 = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn" -> new row 4
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/87e2101526f2d38fc04e67819f1fc856141ab27e/e2e/eb931bac-a78a-4017-93c9-acbdf45de9c6.md",
    "",
    "",
    "eb931bac-a78a-4017-93c9-acbdf45de9c6.md"
) | Out-Null

$wsZhCn.Range("B4").Value2 = "Handed back: in sync with en-US"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c72a545fd95824d4b9964c2b16ac840d28b99fe1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/eb931bac-a78a-4017-93c9-acbdf45de9c6.c72a545fd95824d4b9964c2b16ac840d28b99fe1.zh-cn.xlf",
    "",
    "",
    "eb931bac-a78a-4017-93c9-acbdf45de9c6.c72a545fd95824d4b9964c2b16ac840d28b99fe1.zh-cn.xlf"
) | Out-Null

$wsZhCn.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("D4").Value2 = "2016-02-16 09:40:53"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/159cdeafbb095e158171d85eff3fe9fe076cc94c/e2e/eb931bac-a78a-4017-93c9-acbdf45de9c6.md",
    "",
    "",
    "eb931bac-a78a-4017-93c9-acbdf45de9c6.md"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/36bf0849621298cc9cc261bede42890e22f3fa7d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/eb931bac-a78a-4017-93c9-acbdf45de9c6.c72a545fd95824d4b9964c2b16ac840d28b99fe1.zh-cn.xlf",
    "",
    "",
    "eb931bac-a78a-4017-93c9-acbdf45de9c6.c72a545fd95824d4b9964c2b16ac840d28b99fe1.zh-cn.xlf"
) | Out-Null

$wsZhCn.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("G4").Value2 = "2016-02-16 09:41:45"

$wsZhCn.Range("H4").Value2 = "Include"

# ---------------------------------------------------------------------------
# Sheet 3: "de-de" -> new row 4
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/87e2101526f2d38fc04e67819f1fc856141ab27e/e2e/eb931bac-a78a-4017-93c9-acbdf45de9c6.md",
    "",
    "",
    "eb931bac-a78a-4017-93c9-acbdf45de9c6.md"
) | Out-Null

$wsDeDe.Range("B4").Value2 = "Handed back: in sync with en-US"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c72a545fd95824d4b9964c2b16ac840d28b99fe1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/eb931bac-a78a-4017-93c9-acbdf45de9c6.c72a545fd95824d4b9964c2b16ac840d28b99fe1.de-de.xlf",
    "",
    "",
    "eb931bac-a78a-4017-93c9-acbdf45de9c6.c72a545fd95824d4b9964c2b16ac840d28b99fe1.de-de.xlf"
) | Out-Null

$wsDeDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("D4").Value2 = "2016-02-16 09:41:06"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/25abf46755c98271ec4d2bdf7fdf9ed14c59b995/e2e/eb931bac-a78a-4017-93c9-acbdf45de9c6.md",
    "",
    "",
    "eb931bac-a78a-4017-93c9-acbdf45de9c6.md"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ec9cb55da50e47bd9994a6fd935ff1fb633c66d2/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/eb931bac-a78a-4017-93c9-acbdf45de9c6.c72a545fd95824d4b9964c2b16ac840d28b99fe1.de-de.xlf",
    "",
    "",
    "eb931bac-a78a-4017-93c9-acbdf45de9c6.c72a545fd95824d4b9964c2b16ac840d28b99fe1.de-de.xlf"
) | Out-Null

$wsDeDe.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("G4").Value2 = "2016-02-16 09:42:12"

$wsDeDe.Range("H4").Value2 = "Include"

Write-Host "Generated report row for eb931bac-a78a-4017-93c9-acbdf45de9c6 on Overview, zh-cn and de-de sheets."
